# Redid visualizer code to fit into current flow
#
# Changes applied:
#  1. G2 updated from 30.71099999999998 to 29.81099999999998
#  2. The "Bundle 2" bundle-header block (columns A,B,C,F,G,H) that used
#     to start on row 34 now starts on row 33 instead, with refreshed
#     bundle diameter/weight figures (7.200000000000001 / 6.88).
#     Row 34's corresponding cells are cleared out.
#  3. The vertical merge ranges for columns A,B,C,F,G,H,I,J are shifted
#     accordingly: the "Bundle 1" merge block now ends at row 32 (was 33)
#     and the "Bundle 2" merge block now starts at row 33 (was 34).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unmerge the affected column ranges before touching their cells ---
$ws.Range("A2:A33").UnMerge()
$ws.Range("B2:B33").UnMerge()
$ws.Range("C2:C33").UnMerge()
$ws.Range("F2:F33").UnMerge()
$ws.Range("G2:G33").UnMerge()
$ws.Range("H2:H33").UnMerge()
$ws.Range("I2:I33").UnMerge()
$ws.Range("J2:J33").UnMerge()
$ws.Range("A34:A40").UnMerge()
$ws.Range("B34:B40").UnMerge()
$ws.Range("C34:C40").UnMerge()
$ws.Range("F34:F40").UnMerge()
$ws.Range("G34:G40").UnMerge()
$ws.Range("H34:H40").UnMerge()
$ws.Range("I34:I40").UnMerge()
$ws.Range("J34:J40").UnMerge()

# --- G2 value change ---
$ws.Range("G2").Value = 29.81099999999998

# --- Move the "Bundle 2" header block up from row 34 to row 33 ---
$ws.Range("A33").Value = "Bundle 2"
$ws.Range("B33").Value = "SWITCH-HTR-A"
$ws.Range("C33").Value = "RELAY-RM-1"
$ws.Range("F33").Value = "EXPRESS"
$ws.Range("G33").Value = 7.200000000000001
$ws.Range("H33").Value = 6.88

# --- Clear the cells that used to hold that block on row 34 ---
$ws.Range("A34").Value = ""
$ws.Range("B34").Value = ""
$ws.Range("C34").Value = ""
$ws.Range("F34").Value = ""
$ws.Range("G34").Value = ""
$ws.Range("H34").Value = ""

# --- Re-merge with the updated (shifted by one row) boundaries ---
$ws.Range("A2:A32").Merge()
$ws.Range("B2:B32").Merge()
$ws.Range("C2:C32").Merge()
$ws.Range("F2:F32").Merge()
$ws.Range("G2:G32").Merge()
$ws.Range("H2:H32").Merge()
$ws.Range("I2:I32").Merge()
$ws.Range("J2:J32").Merge()
$ws.Range("A33:A40").Merge()
$ws.Range("B33:B40").Merge()
$ws.Range("C33:C40").Merge()
$ws.Range("F33:F40").Merge()
$ws.Range("G33:G40").Merge()
$ws.Range("H33:H40").Merge()
$ws.Range("I33:I40").Merge()
$ws.Range("J33:J40").Merge()
